$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 37

# Helper: write a value as plain text (shared string), even when the
# text looks like a number/boolean/date, and strip any stray
# number-format / quote-prefix styling that Excel may apply along the
# way so the cell ends up with no explicit style, matching the rest of
# the sheet.
function Set-TextCell($targetRow, $col, $value) {
    $c = $ws.Cells.Item($targetRow, $col)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

# The hyperlink target / screenshot path for the new test run. Add the
# hyperlink (and its display text) first, matching the order in which
# the original authoring tool wrote the shared strings table (the
# screenshot path ends up earlier in sharedStrings.xml than the other
# new strings on this row).
$linkAddress = "C:/Users/AvoComp13/Documents/dlp-automation/test-output/screenshots/2016-09-07-11-26-33-289-AEST.png"
$displayText = "C:\Users\AvoComp13\Documents\dlp-automation\test-output\screenshots\2016-09-07-11-26-33-289-AEST.png"

$apCell = $ws.Cells.Item($row, 42)  # column AP - screenshot
$ws.Hyperlinks.Add($apCell, $linkAddress) | Out-Null
$apCell.Value = $displayText
$apCell.Style = "Normal"

# Fill in the rest of the new row. Most values mirror row 2 (the first
# data row for this device), except for the test name, the two
# "last used" timestamps and the time column, which reflect this new
# test execution.
Set-TextCell $row 1  "Samsung-Galaxy Note5-Generic"            # template
Set-TextCell $row 2  "DEFAULT"                                  # offlineCharging
Set-TextCell $row 3  "Service_NSW], powerControl, reboot, powerSupply, [source, SERVER, offlineCharging, DEFAULT" # roles
Set-TextCell $row 4  "2560"                                     # resolutionHeight
Set-TextCell $row 5  "1"                                        # openRequestsQueue
Set-TextCell $row 6  "English"                                  # language
Set-TextCell $row 7  "1440x2560"                                # resolution
Set-TextCell $row 8  "1115FB5A67CC3205"                         # deviceId
Set-TextCell $row 9  "ec:9b:f3:f1:5d:cd"                        # wifiMacAddress
Set-TextCell $row 10 "Samsung"                                  # manufacturer
Set-TextCell $row 11 "Samsung-Galaxy Note5"                     # templateBaseName
Set-TextCell $row 12 "OPENED"                                   # openingStatus
Set-TextCell $row 13 "SYD-L15O2-14/VIRTUAL/01"                  # cradleId
Set-TextCell $row 14 "deepthi.singh@avocadoconsulting.com.au"   # lockedBy
Set-TextCell $row 15 "5.1.1"                                    # osVersion
Set-TextCell $row 16 "0.01"                                     # testCycle
Set-TextCell $row 17 "params"                                   # testParameters
Set-TextCell $row 18 "true"                                     # inUse
Set-TextCell $row 19 "deepthi.singh@avocadoconsulting.com.au"   # lastUsedBy
Set-TextCell $row 20 "Galaxy Note5"                             # model
Set-TextCell $row 21 "portrait"                                 # id
Set-TextCell $row 22 "nobleltedv-user 5.1.1 LMY47X N920IDVU2AOJ4 release-keys" # firmware
Set-TextCell $row 23 "BrowserTest"                               # testName
Set-TextCell $row 24 "deepthi.singh@avocadoconsulting.com.au"   # allocation
Set-TextCell $row 25 "Android"                                  # os
Set-TextCell $row 26 "rotate"                                   # method
Set-TextCell $row 27 "05.07.2015"                                # releaseDate
Set-TextCell $row 28 "0"                                        # rotation
Set-TextCell $row 29 "Generic"                                  # distributor
Set-TextCell $row 30 "2016-09-07:01-25-19"                      # lastUsedAtFormatted
Set-TextCell $row 31 "SYD-L15O2-14/VIRTUAL/01"                  # lastCradleId
Set-TextCell $row 32 "1473211519267"                            # lastUsedAt
Set-TextCell $row 33 "1115FB5A67CC3205"                         # imei
Set-TextCell $row 34 "APAC-AUS-SYD"                              # location
Set-TextCell $row 35 "2016-09-07 11:26:33"                      # time
Set-TextCell $row 36 "1440"                                     # resolutionWidth
Set-TextCell $row 37 "fail"                                     # testResult
Set-TextCell $row 38 "100"                                      # operabilityRatingScore
Set-TextCell $row 39 "354608070601291"                          # nativeImei
Set-TextCell $row 40 "CONNECTED"                                # status
Set-TextCell $row 41 "mobile"                                   # resourceType
# column 42 (AP, screenshot) was already set above via the hyperlink.
